# Refresh the auto-updating "datetimeFigureOut" date placeholders (footer
# area of the Slide Master and every Slide Layout) so they show the date
# this deck was last edited/saved on, exactly like PowerPoint re-stamps
# those fields whenever the file is saved.
$p = $ppt.ActivePresentation
$newDate = "6/11/2022"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = $null
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ($phType -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    Update-DatePlaceholder $cl.Shapes
}
